# Refresh the "cryptos" price table (rows 2-51) with the latest scrape.
# Price (column D) and Volume/1h (column E) are stored as plain text in this
# sheet (values such as "26.488.88" or "1.838.17" use dots as both thousands
# and decimal separators, so they can never be real numbers) - a leading
# single-quote forces Excel to keep any value that *would* parse as a number
# stored as text, matching the original cell typing.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = '26.488.88'
$ws.Range("E2").Value2 = '  -0.29%  '

# Row 3
$ws.Range("D3").Value2 = '1.838.17'
$ws.Range("E3").Value2 = '  -0.70%  '

# Row 4
$ws.Range("D4").Value2 = '''1.001'
$ws.Range("E4").Value2 = '  -0.04%  '

# Row 5
$ws.Range("D5").Value2 = '''260.63'
$ws.Range("E5").Value2 = '  -0.89%  '

# Row 6
$ws.Range("D6").Value2 = '''1.001'
$ws.Range("E6").Value2 = '  +0.01%  '

# Row 7
$ws.Range("D7").Value2 = '''0.5379'
$ws.Range("E7").Value2 = '  +2.35%  '

# Row 8
$ws.Range("D8").Value2 = '''0.2925'
$ws.Range("E8").Value2 = '  -9.63%  '

# Row 9
$ws.Range("D9").Value2 = '''0.06924'
$ws.Range("E9").Value2 = '  +1.68%  '

# Row 10
$ws.Range("D10").Value2 = '''17.22'
$ws.Range("E10").Value2 = '  -9.04%  '

# Row 11
$ws.Range("D11").Value2 = '1.843.90'
$ws.Range("E11").Value2 = '  -0.22%  '

# Row 12
$ws.Range("D12").Value2 = '''0.7255'
$ws.Range("E12").Value2 = '  -7.52%  '

# Row 13
$ws.Range("D13").Value2 = '''0.07187'
$ws.Range("E13").Value2 = '  -7.50%  '

# Row 14
$ws.Range("D14").Value2 = '''89.04'
$ws.Range("E14").Value2 = '  +0.40%  '

# Row 15
$ws.Range("D15").Value2 = '''4.978'
$ws.Range("E15").Value2 = '  -1.19%  '

# Row 16
$ws.Range("D16").Value2 = '''1.001'
$ws.Range("E16").Value2 = '  +0.04%  '

# Row 17
$ws.Range("D17").Value2 = '''13.75'
$ws.Range("E17").Value2 = '  -1.65%  '

# Row 18
$ws.Range("D18").Value2 = '''1.001'
$ws.Range("E18").Value2 = '  -0.01%  '

# Row 19
$ws.Range("D19").Value2 = '''0.000007868'
$ws.Range("E19").Value2 = '  -1.21%  '

# Row 20
$ws.Range("D20").Value2 = '26.514.89'
$ws.Range("E20").Value2 = '  -0.30%  '

# Row 21
$ws.Range("D21").Value2 = '2.082.92'
$ws.Range("E21").Value2 = '  -0.13%  '

# Row 22
$ws.Range("D22").Value2 = '''4.580'
$ws.Range("E22").Value2 = '  -1.35%  '

# Row 23
$ws.Range("D23").Value2 = '''5.982'
$ws.Range("E23").Value2 = '  -0.34%  '

# Row 24
$ws.Range("D24").Value2 = '''9.179'
$ws.Range("E24").Value2 = '  -3.19%  '

# Row 25
$ws.Range("D25").Value2 = '''141.76'
$ws.Range("E25").Value2 = '  -1.36%  '

# Row 26
$ws.Range("D26").Value2 = '''2.162'
$ws.Range("E26").Value2 = '  -0.55%  '

# Row 27
$ws.Range("D27").Value2 = '''1.704'
$ws.Range("E27").Value2 = '  +1.45%  '

# Row 28
$ws.Range("D28").Value2 = '''16.91'
$ws.Range("E28").Value2 = '  -0.67%  '

# Row 29
$ws.Range("D29").Value2 = '''110.92'
$ws.Range("E29").Value2 = '  -0.94%  '

# Row 30
$ws.Range("D30").Value2 = '''4.225'
$ws.Range("E30").Value2 = '  +0.91%  '

# Row 31
$ws.Range("D31").Value2 = '''0.08877'
$ws.Range("E31").Value2 = '  +1.76%  '

# Row 32
$ws.Range("D32").Value2 = '''4.018'
$ws.Range("E32").Value2 = '  -2.09%  '

# Row 33
$ws.Range("D33").Value2 = '''0.04836'
$ws.Range("E33").Value2 = '  -0.72%  '

# Row 34
$ws.Range("D34").Value2 = '''2.915'
$ws.Range("E34").Value2 = '  +1.39%  '

# Row 35
$ws.Range("D35").Value2 = '''0.7216'
$ws.Range("E35").Value2 = '  -0.16%  '

# Row 36
$ws.Range("D36").Value2 = '''1.128'
$ws.Range("E36").Value2 = '  -0.47%  '

# Row 37
$ws.Range("D37").Value2 = '''3.094'
$ws.Range("E37").Value2 = '  -0.47%  '

# Row 38
$ws.Range("D38").Value2 = '''2.292'
$ws.Range("E38").Value2 = '  +0.78%  '

# Row 39
$ws.Range("D39").Value2 = '''0.01710'
$ws.Range("E39").Value2 = '  -4.57%  '

# Row 40
$ws.Range("D40").Value2 = '''0.4671'
$ws.Range("E40").Value2 = '  -3.82%  '

# Row 41
$ws.Range("D41").Value2 = '''0.9024'
$ws.Range("E41").Value2 = '  +0.17%  '

# Row 42
$ws.Range("D42").Value2 = '''106.87'
$ws.Range("E42").Value2 = '  -3.67%  '

# Row 43
$ws.Range("D43").Value2 = '''5.856'
$ws.Range("E43").Value2 = '  -1.66%  '

# Row 44
$ws.Range("E44").Value2 = '  -0.05%  '

# Row 45
$ws.Range("D45").Value2 = '''7.386'
$ws.Range("E45").Value2 = '  -3.91%  '

# Row 46
$ws.Range("D46").Value2 = '''0.1244'
$ws.Range("E46").Value2 = '  +0.58%  '

# Row 47
$ws.Range("D47").Value2 = '''8.976'
$ws.Range("E47").Value2 = '  -0.92%  '

# Row 48
$ws.Range("B48").Value2 = 'Elrond'
$ws.Range("C48").Value2 = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D48").Value2 = '''34.73'
$ws.Range("E48").Value2 = '  -1.17%  '

# Row 49
$ws.Range("B49").Value2 = 'Decentraland'
$ws.Range("C49").Value2 = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D49").Value2 = '''0.4041'
$ws.Range("E49").Value2 = '  -3.49%  '

# Row 50
$ws.Range("B50").Value2 = 'Cronos'
$ws.Range("C50").Value2 = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value2 = '''0.05755'
$ws.Range("E50").Value2 = '  -2.17%  '

# Row 51
$ws.Range("B51").Value2 = 'EOS'
$ws.Range("C51").Value2 = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D51").Value2 = '''0.8900'
$ws.Range("E51").Value2 = '  -0.26%  '
